# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Agrícola del Norte S.A. de Arica - Palta"
# at row 165, pushing the existing rows 165-176 down to 166-177.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 165 (shifts rows 165:176 down to 166:177)
$ws.Rows.Item(165).Insert()

# Fill in the new row 165 with the latest weekly record
$ws.Cells.Item(165, 1).Value = 1
$ws.Cells.Item(165, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(165, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(165, 4).Value = (Get-Date -Year 2023 -Month 4 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(165, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(165, 5).Value = 15
$ws.Cells.Item(165, 6).Value = "Fruta"
$ws.Cells.Item(165, 7).Value = 100106
$ws.Cells.Item(165, 8).Value = "Oleaginosos"
$ws.Cells.Item(165, 9).Value = 100106002
$ws.Cells.Item(165, 10).Value = "Palta"
$ws.Cells.Item(165, 11).Value = "Hass"
$ws.Cells.Item(165, 12).Value = "Segunda"
$ws.Cells.Item(165, 13).Value = 500
$ws.Cells.Item(165, 14).Value = 30000
$ws.Cells.Item(165, 15).Value = 31000
$ws.Cells.Item(165, 16).Value = 30700
$ws.Cells.Item(165, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(165, 18).Value = "Perú"
$ws.Cells.Item(165, 19).Value = 3070
$ws.Cells.Item(165, 20).Value = 10
